$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 175, shifting existing rows 175-295 down to 176-296.
$ws.Rows.Item(175).Insert()

# Populate the newly inserted row 175 with its data (same categorical
# columns as the rest of the sheet, plus the new date/price figures).
$ws.Range("A175").Value = 8
$ws.Range("B175").Value = "Terminal La Palmera de La Serena"
$ws.Range("C175").Value = "Coquimbo"
$ws.Range("D175").Value = 44981
$ws.Range("E175").Value = 4
$ws.Range("F175").Value = 100112037
$ws.Range("G175").Value = "Cebollín"
$ws.Range("H175").Value = "Sin especificar"
$ws.Range("I175").Value = "Primera"
$ws.Range("J175").Value = 1100
$ws.Range("K175").Value = 1200
$ws.Range("L175").Value = 1400
$ws.Range("M175").Value = 1300
$ws.Range("N175").Value = "$/paquete 6 unidades"
$ws.Range("O175").Value = "Provincia del Elquí"
$ws.Range("P175").Value = 217
$ws.Range("Q175").Value = 6
$ws.Range("R175").Value = "Hortaliza"
